$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.206.51'
$ws.Range("E2").Value = '  +1.06%  '
$ws.Range("D3").Value = '3.341.83'
$ws.Range("E3").Value = '  +0.52%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '582.64'
$ws.Range("E5").Value = '  +0.38%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '176.98'
$ws.Range("E6").Value = '  +1.70%  '
$ws.Range("E7").Value = '  +0.12%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.591'
$ws.Range("E8").Value = '  +0.76%  '
$ws.Range("E9").Value = '  +3.43%  '
$ws.Range("E10").Value = '  +1.30%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '47.98'
$ws.Range("E11").Value = '  +5.76%  '
$ws.Range("E12").Value = '  +1.61%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '691.90'
$ws.Range("E13").Value = '  +4.21%  '
$ws.Range("D14").Value = '3.883.76'
$ws.Range("E14").Value = '  +0.44%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.41'
$ws.Range("E15").Value = '  +0.40%  '
$ws.Range("D16").Value = '68.223.99'
$ws.Range("E16").Value = '  +0.80%  '
$ws.Range("E17").Value = '  +1.34%  '
$ws.Range("D18").Value = '3.341.71'
$ws.Range("E18").Value = '  +0.66%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.44'
$ws.Range("E19").Value = '  +0.08%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.16'
$ws.Range("E20").Value = '  +2.23%  '
$ws.Range("E21").Value = '  +0.83%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.47'
$ws.Range("E22").Value = '  +0.97%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '16.95'
$ws.Range("E23").Value = '  -0.22%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '99.94'
$ws.Range("E24").Value = '  +1.33%  '
$ws.Range("E25").Value = '  +1.94%  '
$ws.Range("E26").Value = '  +0.92%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.51'
$ws.Range("E27").Value = '  +2.98%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '32.98'
$ws.Range("E28").Value = '  -2.21%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.48'
$ws.Range("E29").Value = '  +0.98%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.92'
$ws.Range("E30").Value = '  -6.24%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '558.32'
$ws.Range("E31").Value = '  -5.64%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '11.03'
$ws.Range("E32").Value = '  +1.08%  '
$ws.Range("E33").Value = '  +1.36%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '57.61'
$ws.Range("E34").Value = '  +1.46%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  -0.05%  '
$ws.Range("D36").Value = '3.693.65'
$ws.Range("E36").Value = '  -0.13%  '
$ws.Range("E37").Value = '  +0.92%  '
$ws.Range("E38").Value = '  +3.54%  '
$ws.Range("E39").Value = '  +4.04%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.17'
$ws.Range("E40").Value = '  +2.01%  '
$ws.Range("E41").Value = '  -0.45%  '
$ws.Range("D42").Value = '0.0₃0671'
$ws.Range("E42").Value = '  +1.48%  '
$ws.Range("E43").Value = '  +0.74%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.26'
$ws.Range("E44").Value = '  +1.01%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0412'
$ws.Range("E45").Value = '  +1.67%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.64'
$ws.Range("E46").Value = '  +2.28%  '
$ws.Range("E47").Value = '  +0.82%  '
$ws.Range("E48").Value = '  -0.14%  '
$ws.Range("E49").Value = '  -0.42%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '130.98'
$ws.Range("E50").Value = '  +3.07%  '
$ws.Range("E51").Value = '  +0.44%  '
